$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "'DS1, DS2"
$ws.Range("C3").Value = "'R1, R2"

$ws.Range("F2").Value = 2
$ws.Range("F3").Value = 2
